$d = $word.ActiveDocument

function SplitAt($pos, $markName) {
    # Forces a run boundary at $pos by materialising then immediately
    # discarding a temporary bookmark there (Word always breaks a run
    # around a bookmark anchor, and removing the bookmark leaves the
    # break behind without touching run formatting).
    $d.Bookmarks.Add($markName, $d.Range($pos, $pos)) | Out-Null
    $d.Bookmarks($markName).Delete()
}

# ---------------------------------------------------------------------
# Paragraph: "Mike is 21 and is a full time student. ... class timetable."
# ---------------------------------------------------------------------

# --- Insertion 1: "Mike" + " is very organised and" + " often tries..." ---
$r1 = $d.Content
$r1.Find.Execute("station/university each day. Mike") | Out-Null
$pos1 = $r1.End
$d.Range($pos1, $pos1).InsertAfter(" is very organised and")
$pos1b = $pos1 + (" is very organised and").Length
SplitAt $pos1b "Tmp1b"
SplitAt $pos1 "Tmp1a"

# --- Split mid-word: "...train befor" | "e the one..." ------------------
$r2 = $d.Content
$r2.Find.Execute("often tries to catch the train befor") | Out-Null
$pos2 = $r2.End
SplitAt $pos2 "Tmp2"

# Move the (hidden) _GoBack bookmark to this exact mid-word location,
# matching the target XML's bookmarkStart/bookmarkEnd placement.
$d.Bookmarks.Add("_GoBack", $d.Range($pos2, $pos2)) | Out-Null

# --- Insertion 2: after "cancellations." insert the new sentence -------
$r3 = $d.Content
$r3.Find.Execute("this accounts for any train delays or cancellations.") | Out-Null
$pos3 = $r3.End
$newSentence = " Mike likes to get ahead in his classes and spends his time travelling completing class work and reading lecture notes. "
$d.Range($pos3, $pos3).InsertAfter($newSentence)
$pos3b = $pos3 + $newSentence.Length
SplitAt $pos3b "Tmp3b"
SplitAt $pos3 "Tmp3a"

# ---------------------------------------------------------------------
# "Key Attributes" bullet list: add w:lastRenderedPageBreak before the
# "Parents work" run's text.
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Parents work") | Out-Null
$pos4 = $r4.Start
$lrpb = $d.Range($pos4, $pos4)
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lrpb.InsertXML($xml)
